$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Move the checkmarks from row 14 (Fixed file) down to row 15 (Auto)
# in the "Groups Testing" table: clear B14:E14 and fill B15:F15.

$checkMark = "$([char]0x2714)$([char]0xFE0F)"

$ws.Range("B14:E14").ClearContents()

$ws.Range("B15").Value = $checkMark
$ws.Range("C15").Value = $checkMark
$ws.Range("D15").Value = $checkMark
$ws.Range("E15").Value = $checkMark
$ws.Range("F15").Value = $checkMark

# Update the active selection to match the saved workbook state.
$ws.Range("F15").Select()
